$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the price/volume columns being updated so that
# numeric-looking strings (e.g. "300.85") and percentages (e.g. "-4.76%")
# are kept as literal text, matching the original inline-string cells.
$fmtRange = $ws.Range("D2:E51")
$fmtRange.NumberFormat = "@"

$ws.Range("D2").Value = '300.85'
$ws.Range("E2").Value = '-4.76%'
$ws.Range("D3").Value = '35.17'
$ws.Range("E3").Value = '-1.17%'
$ws.Range("D4").Value = '5.032'
$ws.Range("D5").Value = '0.07957'
$ws.Range("E5").Value = '-1.83%'
$ws.Range("D6").Value = '1.909'
$ws.Range("E6").Value = '-11.19%'
$ws.Range("D7").Value = '7.793'
$ws.Range("E7").Value = '-2.76%'
$ws.Range("D8").Value = '2.929'
$ws.Range("E8").Value = '7.08%'
$ws.Range("D9").Value = '0.9208'
$ws.Range("E9").Value = '-0.68%'
$ws.Range("D10").Value = '0.1354'
$ws.Range("E10").Value = '33.03%'
$ws.Range("D11").Value = '0.1844'
$ws.Range("E11").Value = '-1.42%'
$ws.Range("D12").Value = '0.09566'
$ws.Range("E12").Value = '3.83%'
$ws.Range("D13").Value = '0.03618'
$ws.Range("E13").Value = '0.47%'
$ws.Range("D14").Value = '0.09842'
$ws.Range("D15").Value = '0.001424'
$ws.Range("E15").Value = '-1.49%'
$ws.Range("D16").Value = '0.005786'
$ws.Range("E16").Value = '0.45%'
$ws.Range("D17").Value = '3.511'
$ws.Range("E17").Value = '1.43%'
$ws.Range("D18").Value = '4.037'
$ws.Range("E18").Value = '-2.63%'
$ws.Range("D19").Value = '0.3426'
$ws.Range("E19").Value = '1.80%'
$ws.Range("D20").Value = '0.1311'
$ws.Range("E20").Value = '-1.65%'
$ws.Range("D21").Value = '5.054'
$ws.Range("E21").Value = '-1.61%'
$ws.Range("D22").Value = '0.2465'
$ws.Range("E22").Value = '10.89%'
$ws.Range("D23").Value = '0.04504'
$ws.Range("E23").Value = '-1.71%'
$ws.Range("D24").Value = '0.001217'
$ws.Range("E24").Value = '-2.47%'
$ws.Range("D25").Value = '0.004797'
$ws.Range("E25").Value = '1.94%'
$ws.Range("D26").Value = '0.0001253'
$ws.Range("E26").Value = '0.14%'
$ws.Range("D27").Value = '0.0003008'
$ws.Range("E27").Value = '-33.22%'
$ws.Range("D39").Value = '0.01865'
$ws.Range("E39").Value = '-4.31%'
$ws.Range("D40").Value = '0.04698'
$ws.Range("E40").Value = '-3.40%'
$ws.Range("D41").Value = '0.007554'
$ws.Range("E41").Value = '-2.64%'
$ws.Range("D42").Value = '0.009707'
$ws.Range("E42").Value = '23.94%'
$ws.Range("E43").Value = '-4.59%'
$ws.Range("D44").Value = '0.002115'
$ws.Range("E44").Value = '0.49%'
$ws.Range("D45").Value = '0.01079'
$ws.Range("E45").Value = '-7.18%'
$ws.Range("D46").Value = '0.00006194'
$ws.Range("E46").Value = '-4.64%'
$ws.Range("E47").Value = '0.13%'
$ws.Range("E48").Value = '64.58%'
$ws.Range("E49").Value = '-12.46%'
$ws.Range("D50").Value = '0.00002105'
$ws.Range("E50").Value = '0.13%'
$ws.Range("D51").Value = '0.0002004'
$ws.Range("E51").Value = '0.13%'

# Restore the default cell style now that the values are stored as text,
# so no stray formatting is left behind on the edited cells.
$fmtRange.Style = "Normal"
